$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 33 (pushes old row 33 down to row 35).
$ws.Rows.Item(33).Resize(2).Insert()

# Copy the date style (numeric/date formatting) from the row now at 35 down
# into the two newly-inserted rows so the Fecha column keeps its format.
$ws.Range("D35").Copy()
$ws.Range("D33:D34").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row 33: Angeleno / Segunda
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 45008
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100103
$ws.Range("H33").Value = "Frutos de hueso (carozo)"
$ws.Range("I33").Value = 100103002
$ws.Range("J33").Value = "Ciruela"
$ws.Range("K33").Value = "Angeleno"
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 300
$ws.Range("N33").Value = 19000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 19500
$ws.Range("Q33").Value = "$/bandeja 18 kilos granel"
$ws.Range("R33").Value = "Región de O'Higgins"
$ws.Range("S33").Value = 1083
$ws.Range("T33").Value = 18

# New row 34: Fortuna / Tercera
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 45008
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103002
$ws.Range("J34").Value = "Ciruela"
$ws.Range("K34").Value = "Fortuna"
$ws.Range("L34").Value = "Tercera"
$ws.Range("M34").Value = 270
$ws.Range("N34").Value = 19000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 19500
$ws.Range("Q34").Value = "$/bandeja 18 kilos granel"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 1083
$ws.Range("T34").Value = 18
